# Guardar Excel en Datagridview
# Add a new column (D) to Hoja1 and create a new sheet (Hoja2) that
# lists the column headers of Hoja1 together with the row values of
# the new column D.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# New column D on Hoja1: header "D", values 1..4
$ws1.Range("D1").Value = "D"
$ws1.Range("D2").Value = 1
$ws1.Range("D3").Value = 2
$ws1.Range("D4").Value = 3
$ws1.Range("D5").Value = 4

# Add the new worksheet, placed right after Hoja1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

# Header row: "D", "E"
$ws2.Range("A1").Value = "D"
$ws2.Range("B1").Value = "E"

# Data rows: numbers 1..4 paired with original headers A, B, C, D
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "A"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "B"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "C"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "D"

$ws2.Range("B5").Select()

# Re-activate Hoja1 and restore its selection so it remains the visible tab
$ws1.Activate()
$ws1.Range("D1:D5").Select()
